$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3: swap Team name and Weekly Pending Total(Rp), update Repayment
$ws.Range("B2").Value = "Xinghao_s2l"
$ws.Range("C2").Value = 2540164470
$ws.Range("D2").Value = 277031612

$ws.Range("B3").Value = "Cpu_s2l"
$ws.Range("C3").Value = 1705546967
$ws.Range("D3").Value = 184313825

# Rows 4-8: updated Repayment figures only
$ws.Range("D4").Value = 663367261
$ws.Range("D5").Value = 171667734
$ws.Range("D6").Value = 594832964
$ws.Range("D7").Value = 561131315
$ws.Range("D8").Value = 300291681

# Update the selection to match the saved view state
$ws.Range("C2").Select()
